$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.991.06'
$ws.Range('E2').Value = '  +1.19%  '
$ws.Range('D3').Value = '2.508.89'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''323.44'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('D6').Value = '''108.46'
$ws.Range('E6').Value = '  -0.99%  '
$ws.Range('D7').Value = '''0.524'
$ws.Range('E7').Value = '  -0.66%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '''0.561'
$ws.Range('E9').Value = '  +3.48%  '
$ws.Range('D10').Value = '''40.27'
$ws.Range('E10').Value = '  +3.12%  '
$ws.Range('D11').Value = '''19.69'
$ws.Range('E11').Value = '  +5.43%  '
$ws.Range('D12').Value = '''0.0815'
$ws.Range('E12').Value = '  -0.53%  '
$ws.Range('E13').Value = '  +0.65%  '
$ws.Range('D14').Value = '''7.19'
$ws.Range('E14').Value = '  -0.35%  '
$ws.Range('D15').Value = '2.898.04'
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('D16').Value = '2.509.45'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').Value = '''0.850'
$ws.Range('E17').Value = '  -1.67%  '
$ws.Range('D18').Value = '47.843.63'
$ws.Range('E18').Value = '  +0.97%  '
$ws.Range('E19').Value = '  +2.29%  '
$ws.Range('D20').Value = '''6.61'
$ws.Range('E20').Value = '  -1.70%  '
$ws.Range('D21').Value = '0.0₃0942'
$ws.Range('E21').Value = '  -0.74%  '
$ws.Range('E22').Value = '  +4.54%  '
$ws.Range('D23').Value = '''70.93'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').Value = '''247.64'
$ws.Range('E24').Value = '  -0.90%  '
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = '''25.86'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').Value = '''10.22'
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('D29').Value = '''2.27'
$ws.Range('E29').Value = '  -1.30%  '
$ws.Range('E30').Value = '  +4.37%  '
$ws.Range('D31').Value = '''35.12'
$ws.Range('E31').Value = '  -2.52%  '
$ws.Range('D32').Value = '''49.82'
$ws.Range('E32').Value = '  -0.86%  '
$ws.Range('D33').Value = '''19.99'
$ws.Range('E33').Value = '  +0.14%  '
$ws.Range('D34').Value = '''5.38'
$ws.Range('E34').Value = '  -1.41%  '
$ws.Range('D35').Value = '''1.00'
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('D36').Value = '''0.0785'
$ws.Range('E36').Value = '  -1.26%  '
$ws.Range('D37').Value = '''1.98'
$ws.Range('E37').Value = '  -1.30%  '
$ws.Range('D38').Value = '''4.69'
$ws.Range('E38').Value = '  -1.64%  '
$ws.Range('E39').Value = '  -1.16%  '
$ws.Range('E40').Value = '  -0.70%  '
$ws.Range('D41').Value = '''22.23'
$ws.Range('E41').Value = '  +3.70%  '
$ws.Range('D42').Value = '''118.90'
$ws.Range('E42').Value = '  -3.14%  '
$ws.Range('E43').Value = '  -3.57%  '
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('D45').Value = '1.996.79'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').Value = '''3.12'
$ws.Range('E46').Value = '  +1.22%  '
$ws.Range('E47').Value = '  -2.90%  '
$ws.Range('E48').Value = '  +0.91%  '
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('D50').Value = '''5.19'
$ws.Range('E50').Value = '  -2.10%  '
$ws.Range('D51').Value = '''56.59'
$ws.Range('E51').Value = '  +0.95%  '
